# Golden Shoe proposal - add "Discount vouchers" backlog item + return
# process paragraph notes. Mirrors the author's commit:
#   "Added support and some backend return functionality"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: relocate the _GoBack bookmark.
#
# In the committed revision the _GoBack bookmark (which Word stamps at
# the location of the most recent edit) ends up right after the
# "Return process" paragraph's text, instead of at the very start of
# the document (its original position here). We reproduce that by
# deleting the existing bookmark and re-adding it in the new spot.
# ---------------------------------------------------------------------

$anchor = $d.Content
$found = $anchor.Find.Execute("for a given item. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $targetPos = $anchor.End

    # Drop a temporary marker right at the target position so we have a
    # stable, non-boundary place to collapse a Range to (collapsed
    # ranges exactly at a paragraph's last-character boundary don't
    # resolve correctly for Bookmarks.Add in this host).
    $markerAnchor = $d.Range($targetPos, $targetPos)
    $markerAnchor.InsertAfter("@@GOBACKMARK@@")

    $existingBookmark = $d.Bookmarks.Item("_GoBack")
    $existingBookmark.Delete()

    $newBookmarkRange = $d.Range($targetPos, $targetPos)
    $d.Bookmarks.Add("_GoBack", $newBookmarkRange)

    $markerRange = $d.Content
    $markerRange.Find.Execute("@@GOBACKMARK@@", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $markerRange.Delete()
}

# ---------------------------------------------------------------------
# Part 2: append the new "Discount vouchers" backlog section at the end
# of the document.
# ---------------------------------------------------------------------

# -- Heading2: "Discount vouchers ..." --
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Range.Text = "Discount vouchers that can only be used in store send directly to customer via email"

# -- Body paragraph describing the problem --
$headingPara2 = $d.Paragraphs.Last
$headingPara2.Range.InsertParagraphAfter()
$bodyPara1 = $d.Paragraphs.Last
$bodyPara1.Range.Text = "I believe this is quite inefficient as staff members have to manually (I assume) collect user information and send emails with discount vouchers."

# -- "Solution: ..." paragraph (bold lead-in) --
$bodyPara1b = $d.Paragraphs.Last
$bodyPara1b.Range.InsertParagraphAfter()
$solutionPara = $d.Paragraphs.Last
$solutionPara.Range.Text = "Solution: Add an “Apply voucher” functionality on check-out that will automatically check for existing voucher and apply a discount."

# -- Two trailing blank paragraphs --
$solutionPara2 = $d.Paragraphs.Last
$solutionPara2.Range.InsertParagraphAfter()
$blankPara1 = $d.Paragraphs.Last

$blankPara1.Range.InsertParagraphAfter()
$blankPara2 = $d.Paragraphs.Last

# Now apply the Heading2 style to just the new heading paragraph -
# done after all the siblings exist so the style doesn't get inherited
# forward onto the body/blank paragraphs that follow it.
$headingParaFinal = $d.Paragraphs.Item($headingPara.Index)
$headingParaFinal.Style = "Heading 2"

# Bold the "Solution: " lead-in of the new solution paragraph.
$solutionParaFinal = $d.Paragraphs.Item($solutionPara.Index)
$boldRange = $d.Range($solutionParaFinal.Range.Start, $solutionParaFinal.Range.Start + 10)
$boldRange.Font.Bold = 1

Write-Output "edit applied"
